$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Emission type column (C) was blank/empty numeric for data rows 2-4.
# Correct it to "C" (alternate receptor use correction) for each pollutant row.
$ws.Range("C2").Value = "C"
$ws.Range("C3").Value = "C"
$ws.Range("C4").Value = "C"
